function Set-TextValue($range, $val) {
    # Forces a value to be stored as text, even if it looks like a number
    # (e.g. "6.60" or "008227"), without leaving stray direct formatting behind.
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1 "2021-Q4" is left untouched.
# Sheet 2 was named "总计"; it becomes "2022-Q1" and gets brand-new data.
# A new sheet named "总计" is appended at the end with updated summary data.
# ---------------------------------------------------------------------

$wsQ1 = $wb.Worksheets.Item(1)          # "2021-Q4" -- used as a formatting donor
$wsQ22022 = $wb.Worksheets.Item(2)      # currently "总计"

# Rename the second sheet to "2022-Q1" first so the name "总计" is free to reuse later.
$wsQ22022.Name = "2022-Q1"

# --- Header row for "2022-Q1" (columns B..H), reusing the bold/border style ("s=2") ---
$headers2022 = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers2022.Length; $i++) {
    $col = 2 + $i   # B=2 .. H=8
    $cell = $wsQ22022.Cells.Item(1, $col)
    $cell.Value = $headers2022[$i]
    $wsQ1.Range("B1").Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats, reuse existing style "s=2"
}
$wsQ22022.Application.CutCopyMode = 0

# --- Data rows for "2022-Q1" ---
# Each entry: index(A,numeric), code(B,text), name(C,text), size(D,text),
#             stockPos(E,text), posPct(F,text), marketValue(G,text), rank(H,numeric)
$rows2022 = @(
    @(0, "008227", "宝盈研究精选混合A",       "6.60", "92.88", "3.38", "0.2231", 10),
    @(1, "008228", "宝盈研究精选混合C",       "1.67", "92.88", "3.38", "0.0564", 10),
    @(2, "350001", "天治财富增长混合",         "0.98", "69.00", "3.33", "0.0326", 2),
    @(3, "004890", "中邮健康文娱灵活配置混合", "0.44", "81.43", "5.11", "0.0225", 5)
)

$rowNum = 2
foreach ($r in $rows2022) {
    # Column A: numeric index value, reuse the "s=2" style already present on A2.
    $wsQ1.Range("A2").Copy() | Out-Null
    $wsQ22022.Cells.Item($rowNum, 1).PasteSpecial(-4122) | Out-Null
    $wsQ22022.Cells.Item($rowNum, 1).Value = $r[0]

    # Columns B, C, D, E, F, G: text values (B and the numeric-looking ones must
    # not be auto-coerced into numbers -- e.g. "008227" must keep its leading zero).
    Set-TextValue $wsQ22022.Cells.Item($rowNum, 2) $r[1]
    Set-TextValue $wsQ22022.Cells.Item($rowNum, 3) $r[2]
    Set-TextValue $wsQ22022.Cells.Item($rowNum, 4) $r[3]
    Set-TextValue $wsQ22022.Cells.Item($rowNum, 5) $r[4]
    Set-TextValue $wsQ22022.Cells.Item($rowNum, 6) $r[5]
    Set-TextValue $wsQ22022.Cells.Item($rowNum, 7) $r[6]

    # Column H: real number.
    $wsQ22022.Cells.Item($rowNum, 8).Value = $r[7]

    $rowNum++
}
$wsQ22022.Application.CutCopyMode = 0

# ---------------------------------------------------------------------
# Add the new "总计" sheet at the very end of the workbook. It is created by
# copying an existing sheet (to inherit the same sheetPr/pageMargins/etc
# structure) and then wiping its cells completely before laying out the
# new summary table.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsQ1.Copy($null, $lastSheet) | Out-Null
$wsZongji = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsZongji.Cells.Clear() | Out-Null
$wsZongji.Name = "总计"

# Header row (reuse "s=2" style).
$headersZongji = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($i = 0; $i -lt $headersZongji.Length; $i++) {
    $col = 2 + $i   # B=2 .. D=4
    $cell = $wsZongji.Cells.Item(1, $col)
    $cell.Value = $headersZongji[$i]
    $wsQ1.Range("B1").Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}
$wsZongji.Application.CutCopyMode = 0

# Data rows: (index, date label, count, market value)
$rowsZongji = @(
    @(0, "2022-Q1", 4, 0.33),
    @(1, "2021-Q4", 6, 0.73)
)

$rowNum = 2
foreach ($r in $rowsZongji) {
    $wsQ1.Range("A2").Copy() | Out-Null
    $wsZongji.Cells.Item($rowNum, 1).PasteSpecial(-4122) | Out-Null
    $wsZongji.Cells.Item($rowNum, 1).Value = $r[0]

    Set-TextValue $wsZongji.Cells.Item($rowNum, 2) $r[1]
    $wsZongji.Cells.Item($rowNum, 3).Value = $r[2]
    $wsZongji.Cells.Item($rowNum, 4).Value = $r[3]

    $rowNum++
}
$wsZongji.Application.CutCopyMode = 0
